$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1805.375
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1991.8572
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 5975.571599999999
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -6311.571599999999
# Row 41
$ws.Range("H41").Value = 866
$ws.Range("I41").Value = 1423.4445
$ws.Range("J41").Value = 149.28572
$ws.Range("K41").Value = 1423.4445
$ws.Range("L41").Value = 149.28572
$ws.Range("M41").Value = -983.4445000000001
$ws.Range("N41").Value = -1029.28572
# Row 74
$ws.Range("H74").Value = 10440.454
$ws.Range("I74").Value = 10192.143
$ws.Range("J74").Value = 10875
$ws.Range("K74").Value = 10192.143
$ws.Range("L74").Value = 10875
$ws.Range("M74").Value = -9256.143
$ws.Range("N74").Value = -12747
# Row 77
$ws.Range("H77").Value = 10440.454
$ws.Range("I77").Value = 10192.143
$ws.Range("J77").Value = 10875
$ws.Range("K77").Value = 50960.715
$ws.Range("L77").Value = 54375
$ws.Range("M77").Value = -46280.715
$ws.Range("N77").Value = -63735
# Row 95
$ws.Range("H95").Value = 17688.6
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17688.6
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17688.6
$ws.Range("N95").Value = -23180.6
# Row 132
$ws.Range("H132").Value = 4570.4
$ws.Range("I132").Value = 4570.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13711.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11181.2
$ws.Range("N132").ClearContents()
# Row 137
$ws.Range("H137").Value = 2318.8572
$ws.Range("I137").Value = 2448.4
$ws.Range("J137").Value = 2246.889
$ws.Range("K137").Value = 7345.200000000001
$ws.Range("L137").Value = 6740.667
$ws.Range("M137").Value = -4795.200000000001
$ws.Range("N137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2303.5557
$ws.Range("I2").Value = 2510.875
$ws.Range("J2").Value = 645
$ws.Range("K2").Value = 2510.875
$ws.Range("L2").Value = 645
$ws.Range("M2").Value = -2397.875
$ws.Range("N2").ClearContents()
# Row 32
$ws.Range("H32").Value = 8617.964
$ws.Range("I32").Value = 760.4
$ws.Range("J32").Value = 40762.547
$ws.Range("K32").Value = 760.4
$ws.Range("L32").Value = 40762.547
$ws.Range("M32").Value = -473.4
$ws.Range("N32").ClearContents()
# Row 95
$ws.Range("H95").Value = 31779.8
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 31779.8
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 31779.8
$ws.Range("N95").Value = -37271.8
# Row 96
$ws.Range("H96").Value = 44995
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 44995
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 44995
$ws.Range("N96").Value = -50487
# Row 97
$ws.Range("H97").Value = 569.4
$ws.Range("I97").Value = 608.7273
$ws.Range("J97").Value = 461.25
$ws.Range("K97").Value = 608.7273
$ws.Range("L97").Value = 461.25
$ws.Range("M97").Value = -112.7273
$ws.Range("N97").Value = -1453.25
# Row 116
$ws.Range("H116").Value = 2303.5557
$ws.Range("I116").Value = 2510.875
$ws.Range("J116").Value = 645
$ws.Range("K116").Value = 2510.875
$ws.Range("L116").Value = 645
$ws.Range("M116").Value = -216.875
$ws.Range("N116").ClearContents()
# Row 132
$ws.Range("H132").Value = 2932.032
$ws.Range("I132").Value = 2747.4614
$ws.Range("J132").Value = 8530.666999999999
$ws.Range("K132").Value = 8242.3842
$ws.Range("L132").Value = 25592.001
$ws.Range("M132").Value = -5712.3842
$ws.Range("N132").Value = -30652.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2303.5557
$ws.Range("I3").Value = 2510.875
$ws.Range("J3").Value = 645
$ws.Range("K3").Value = 2510.875
$ws.Range("L3").Value = 645
$ws.Range("M3").Value = -2396.875
$ws.Range("N3").ClearContents()
# Row 105
$ws.Range("H105").Value = 4458.8667
$ws.Range("I105").Value = 4063.0715
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 4063.0715
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -2316.0715
$ws.Range("N105").ClearContents()
# Row 134
$ws.Range("H134").Value = 6137.1333
$ws.Range("I134").Value = 5547.1665
$ws.Range("J134").Value = 8497
$ws.Range("K134").Value = 16641.4995
$ws.Range("L134").Value = 25491
$ws.Range("M134").Value = -14106.4995
$ws.Range("N134").Value = -30561

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 2243
$ws.Range("I5").Value = 235.66667
$ws.Range("J5").Value = 5254
$ws.Range("K5").Value = 235.66667
$ws.Range("L5").Value = 5254
$ws.Range("M5").Value = -123.66667
$ws.Range("N5").Value = -5478
# Row 10
$ws.Range("H10").Value = 5004.6665
$ws.Range("I10").Value = 1007
$ws.Range("J10").Value = 7003.5
$ws.Range("K10").Value = 1007
$ws.Range("L10").Value = 7003.5
$ws.Range("M10").Value = -868
$ws.Range("N10").ClearContents()
# Row 22
$ws.Range("H22").Value = 490
$ws.Range("I22").Value = 475
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 475
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -125
$ws.Range("N22").ClearContents()
# Row 62
$ws.Range("H62").Value = 3192
$ws.Range("I62").Value = 3225.6667
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 3225.6667
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -2601.6667
$ws.Range("N62").Value = -4238
# Row 65
$ws.Range("H65").Value = 3192
$ws.Range("I65").Value = 3225.6667
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 16128.3335
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -13008.3335
$ws.Range("N65").Value = -21190

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 8102.222
$ws.Range("I39").Value = 3795
$ws.Range("J39").Value = 8355.588
$ws.Range("K39").Value = 11385
$ws.Range("L39").Value = 25066.764
$ws.Range("M39").Value = -11091
$ws.Range("N39").Value = -25654.764
# Row 55
$ws.Range("H55").Value = 7162.6294
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 7162.6294
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21487.8882
$ws.Range("N55").Value = -21841.8882
# Row 80
$ws.Range("H80").Value = 3794.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3794.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11383.5
$ws.Range("N80").Value = -13255.5
# Row 83
$ws.Range("H83").Value = 3794.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3794.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 34150.5
$ws.Range("N83").Value = -43510.5
# Row 103
$ws.Range("H103").Value = 629.36365
$ws.Range("I103").Value = 627.75
$ws.Range("J103").Value = 633.6667
$ws.Range("K103").Value = 1883.25
$ws.Range("L103").Value = 1901.0001
$ws.Range("M103").Value = -1004.25
$ws.Range("N103").ClearContents()
# Row 104
$ws.Range("H104").Value = 4235
$ws.Range("I104").Value = 2000
$ws.Range("J104").Value = 4980
$ws.Range("K104").Value = 6000
$ws.Range("L104").Value = 14940
$ws.Range("M104").Value = -3379
$ws.Range("N104").ClearContents()
# Row 118
$ws.Range("H118").Value = 7999.6665
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 7999.6665
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 23998.9995
$ws.Range("N118").Value = -26484.9995
# Row 133
$ws.Range("H133").Value = 11590.5
$ws.Range("I133").Value = 4893.5
$ws.Range("J133").Value = 18287.5
$ws.Range("K133").Value = 14680.5
$ws.Range("L133").Value = 54862.5
$ws.Range("M133").Value = -9620.5
$ws.Range("N133").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 180
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 180
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 180
$ws.Range("N13").Value = -458
# Row 48
$ws.Range("H48").Value = 33353334
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 33353334
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 33353334
$ws.Range("N48").Value = -33354304
# Row 59
$ws.Range("H59").Value = 46599.8
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 50749.75
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 50749.75
$ws.Range("M59").Value = -29417
$ws.Range("N59").Value = -51915.75
# Row 95
$ws.Range("H95").Value = 28950
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 28950
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 28950
$ws.Range("N95").Value = -34442
# Row 102
$ws.Range("H102").Value = 2918.6428
$ws.Range("I102").Value = 2967.077
$ws.Range("J102").Value = 2289
$ws.Range("K102").Value = 2967.077
$ws.Range("L102").Value = 2289
$ws.Range("M102").Value = -1345.077
$ws.Range("N102").Value = -5533
# Row 122
$ws.Range("H122").Value = 1998
$ws.Range("I122").Value = 1998
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5994
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3544
$ws.Range("N122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 17150
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 17150
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 17150
$ws.Range("N64").Value = -17600
# Row 67
$ws.Range("H67").Value = 17150
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 17150
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 17150
$ws.Range("N67").Value = -18710
# Row 104
$ws.Range("H104").Value = 9567.6
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 9567.6
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 9567.6
$ws.Range("N104").Value = -16555.6
# Row 122
$ws.Range("H122").Value = 8765.429
$ws.Range("I122").Value = 9186
$ws.Range("J122").Value = 3298
$ws.Range("K122").Value = 27558
$ws.Range("L122").Value = 9894
$ws.Range("M122").Value = -25108
$ws.Range("N122").Value = -14794

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 2289571.5
$ws.Range("I4").Value = 5168000
$ws.Range("J4").Value = 130750
$ws.Range("K4").Value = 5168000
$ws.Range("L4").Value = 130750
$ws.Range("M4").Value = -5167887
$ws.Range("N4").Value = -130976
# Row 63
$ws.Range("H63").Value = 42000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 42000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 42000
$ws.Range("N63").Value = -43248
# Row 66
$ws.Range("H66").Value = 42000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 42000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 126000
$ws.Range("N66").Value = -132240
# Row 81
$ws.Range("H81").Value = 2638.8
$ws.Range("I81").Value = 2048.5
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 4097
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -3036
$ws.Range("N81").Value = -12122
# Row 84
$ws.Range("H84").Value = 2638.8
$ws.Range("I84").Value = 2048.5
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 20485
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -15181
$ws.Range("N84").Value = -60608
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
# Row 119
$ws.Range("H119").Value = 1576833.2
$ws.Range("I119").Value = 6000000
$ws.Range("J119").Value = 102444.336
$ws.Range("K119").Value = 6000000
$ws.Range("L119").Value = 102444.336
$ws.Range("M119").Value = -5995162
$ws.Range("N119").Value = -112120.336
# Row 136
$ws.Range("H136").Value = 2928.1143
$ws.Range("I136").Value = 2446.7307
$ws.Range("J136").Value = 4318.778
$ws.Range("K136").Value = 7340.1921
$ws.Range("L136").Value = 12956.334
$ws.Range("M136").Value = -4790.1921
$ws.Range("N136").ClearContents()
